# "19.05.19 Today Sales Updated"
# Update the "Raju Ahamed" sheet's date header (18.05.19 -> 19.05.19)
# and today's sales quantities, for both duplicated blocks on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Raju Ahamed")

# Date header cells (both copies on this sheet show the same date string)
$ws.Range("A4").Value = "Date: 19.05.19"
$ws.Range("A31").Value = "Date: 19.05.19"

# Today's sales quantities - first block (rows 6-10)
$ws.Range("E6").Value = 37
$ws.Range("E7").Value = 110
$ws.Range("E8").Value = 91
$ws.Range("E9").Value = 46
$ws.Range("E10").Value = 30

# Today's sales quantities - second (duplicate) block (rows 33-37)
$ws.Range("E33").Value = 37
$ws.Range("E34").Value = 110
$ws.Range("E35").Value = 91
$ws.Range("E36").Value = 46
$ws.Range("E37").Value = 30
